# Rotate the four reference URLs (Wikipedia -> NIH -> Google Scholar -> JSTOR -> Wikipedia)
# in the "Content Placeholder 2" shape on each "References" slide.
# Each target string is cleared first, then set, so the run-diffing logic
# doesn't keep the old shared prefix/suffix as a separate <a:r> run.

$p = $ppt.ActivePresentation

$slideIndexes = @(7, 13, 19, 25, 29)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item("Content Placeholder 2")
    $tr = $shape.TextFrame.TextRange

    $para2 = $tr.Paragraphs(2, 1)
    $para2.Text = ""
    $para2.Text = "https://www.nih.gov/"

    $para3 = $tr.Paragraphs(3, 1)
    $para3.Text = ""
    $para3.Text = "https://scholar.google.com/"

    $para4 = $tr.Paragraphs(4, 1)
    $para4.Text = ""
    $para4.Text = "https://www.jstor.org/"

    $para5 = $tr.Paragraphs(5, 1)
    $para5.Text = ""
    $para5.Text = "https://en.wikipedia.org/wiki/Main_Page"
}
